# Updated symbol list on Tue Dec 13 23:42:35 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text, even though they look
# numeric, so every new value is written with a leading apostrophe to
# force Excel to keep it as a text string instead of re-interpreting it
# as a number (which would introduce floating point noise and change
# the cell type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range("D2").Value  = "'271.85"
$ws.Range("D3").Value  = "'22.88"
$ws.Range("D4").Value  = "'6.481"
$ws.Range("D5").Value  = "'0.06218"
$ws.Range("D6").Value  = "'3.648"
$ws.Range("D7").Value  = "'6.652"
$ws.Range("D9").Value  = "'0.8298"
$ws.Range("D10").Value = "'0.01381"
$ws.Range("D11").Value = "'0.1600"
$ws.Range("D12").Value = "'0.08301"
$ws.Range("D13").Value = "'0.03434"
$ws.Range("D14").Value = "'0.03188"
$ws.Range("D16").Value = "'3.845"
$ws.Range("D17").Value = "'0.001644"
$ws.Range("D19").Value = "'0.006321"
$ws.Range("D20").Value = "'0.005693"
$ws.Range("D21").Value = "'0.001076"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.717"
$ws.Range("D24").Value = "'2.325"
$ws.Range("D27").Value = "'0.0002700"
$ws.Range("D40").Value = "'0.04702"
$ws.Range("D41").Value = "'0.007055"

# --- Row 42 and 43: BKEXToken and CEJI swap places (values are not a
#     straight swap, the underlying data also changed) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1162"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003349"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

$ws.Range("D44").Value = "'0.01204"
$ws.Range("D45").Value = "'0.00006240"
$ws.Range("D48").Value = "'0.9188"
$ws.Range("D49").Value = "'0.002132"

$ws.Range("D50").Value = "'0.00001398"
$ws.Range("E50").Value = "49CryptobidCoinCBC"

$ws.Range("D51").Value = "'0.01238"
